$d = $word.ActiveDocument

# Locate the full sentence that needs to be split into three runs and
# corrected ("possibile" -> "possibili").
$oldFull = " più esteso e decisamente più training, porta risultati promettenti sulle possibile capacità che potrebbe avere senza le mancanze che ha in questo momento."

$full = $d.Content
$found = $full.Find.Execute($oldFull, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Target sentence not found in document."
}

$base = $full.Start

# The three pieces the single run must be split into (text content only;
# formatting - lang=it-IT - stays identical across all three).
$seg1    = " più esteso e decisamente più training, porta risul"
$seg2Old = "tati promettenti sulle possibile"
$seg2New = "tati promettenti sulle possibili"
$seg3    = " capacità che potrebbe avere senza le mancanze che ha in questo momento."

$s1Start = $base
$s1End   = $s1Start + $seg1.Length

$s2Start = $s1End
$s2End   = $s2Start + $seg2Old.Length

$s3Start = $s2End
$s3End   = $s3Start + $seg3.Length

# First fix the typo in-place (this keeps the run merged for now, but the
# text content is correct afterwards).
$r2 = $d.Range($s2Start, $s2End)
$r2.Text = $seg2New

# Now force the surrounding text to split into independent runs by
# toggling a character-formatting property on/off for the first two
# segments; this creates genuine run boundaries without altering the
# final (identical) formatting.
$r1 = $d.Range($s1Start, $s1End)
$r1.Font.Bold = 1
$r1.Font.Bold = 0

$r2b = $d.Range($s2Start, $s2End)
$r2b.Font.Bold = 1
$r2b.Font.Bold = 0
